$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: the match data (cols F:V) for row 37 and row 38 was swapped ---
# (A:E - index/country/tournament/season/date - stay as-is, only F:V change)
# Row 37
$ws.Range("F37").Value = "Al Bataeh"
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = "Ajman"
$ws.Range("I37").Value = 1
$ws.Range("J37").Value = 2.67
$ws.Range("K37").Value = "23/10/2023 17:42"
$ws.Range("L37").Value = 2.36
$ws.Range("M37").Value = "27/10/2023 14:40"
$ws.Range("N37").Value = 3.68
$ws.Range("O37").Value = "23/10/2023 17:42"
$ws.Range("P37").Value = 3.82
$ws.Range("Q37").Value = "27/10/2023 14:36"
$ws.Range("R37").Value = 2.36
$ws.Range("S37").Value = "23/10/2023 17:42"
$ws.Range("T37").Value = 2.79
$ws.Range("U37").Value = "27/10/2023 14:40"
$ws.Range("V37").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-bataeh-ajman/zLLLBlCo/"

# Row 38
$ws.Range("F38").Value = "Hatta"
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = "Al Wasl"
$ws.Range("I38").Value = 5
$ws.Range("J38").Value = 6.3
$ws.Range("K38").Value = "26/10/2023 13:42"
$ws.Range("L38").Value = 10.27
$ws.Range("M38").Value = "27/10/2023 14:13"
$ws.Range("N38").Value = 4.97
$ws.Range("O38").Value = "26/10/2023 13:42"
$ws.Range("P38").Value = 7.25
$ws.Range("Q38").Value = "27/10/2023 14:13"
$ws.Range("R38").Value = 1.38
$ws.Range("S38").Value = "26/10/2023 13:42"
$ws.Range("T38").Value = 1.22
$ws.Range("U38").Value = "27/10/2023 13:29"
$ws.Range("V38").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/hatta-al-wasl/SYQu7QQG/"

# --- Step 2: append 6 new match rows (52-57) ---
# Clone formatting (styles) from row 51 first so the new rows look identical
# (bold/bordered index in col A, datetime format in col E, plain cells elsewhere)
$ws.Range("A51:V51").Copy()
$ws.Range("A52:V57").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 52 (Indice=51)
$ws.Range("B52").Value = "united-arab-emirates"
$ws.Range("C52").Value = "uae-league"
$ws.Range("D52").Value = "2023-2024"
$ws.Range("A52").Value = 51
$ws.Range("E52").Value = 45255.57291666666
$ws.Range("F52").Value = "Al Wahda"
$ws.Range("G52").Value = 4
$ws.Range("H52").Value = "Emirates Club"
$ws.Range("I52").Value = 1
$ws.Range("J52").Value = 1.33
$ws.Range("K52").Value = "24/11/2023 15:12"
$ws.Range("L52").Value = 1.31
$ws.Range("M52").Value = "25/11/2023 13:36"
$ws.Range("N52").Value = 5.27
$ws.Range("O52").Value = "24/11/2023 15:12"
$ws.Range("P52").Value = 5.9
$ws.Range("Q52").Value = "25/11/2023 13:36"
$ws.Range("R52").Value = 8.22
$ws.Range("S52").Value = "24/11/2023 15:12"
$ws.Range("T52").Value = 8.33
$ws.Range("U52").Value = "25/11/2023 13:36"
$ws.Range("V52").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-wahda-emirates-club/63dfCL2d/"

# Row 53 (Indice=52)
$ws.Range("B53").Value = "united-arab-emirates"
$ws.Range("C53").Value = "uae-league"
$ws.Range("D53").Value = "2023-2024"
$ws.Range("A53").Value = 52
$ws.Range("E53").Value = 45255.57291666666
$ws.Range("F53").Value = "Hatta"
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = "Bani Yas"
$ws.Range("I53").Value = 1
$ws.Range("J53").Value = 3.08
$ws.Range("K53").Value = "24/11/2023 15:12"
$ws.Range("L53").Value = 4.35
$ws.Range("M53").Value = "25/11/2023 13:43"
$ws.Range("N53").Value = 3.72
$ws.Range("O53").Value = "24/11/2023 15:12"
$ws.Range("P53").Value = 4.42
$ws.Range("Q53").Value = "25/11/2023 13:43"
$ws.Range("R53").Value = 2.09
$ws.Range("S53").Value = "24/11/2023 15:12"
$ws.Range("T53").Value = 1.69
$ws.Range("U53").Value = "25/11/2023 13:43"
$ws.Range("V53").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/hatta-bani-yas/6uE7LPuA/"

# Row 54 (Indice=53)
$ws.Range("B54").Value = "united-arab-emirates"
$ws.Range("C54").Value = "uae-league"
$ws.Range("D54").Value = "2023-2024"
$ws.Range("A54").Value = 53
$ws.Range("E54").Value = 45255.6875
$ws.Range("F54").Value = "Al Jazira"
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = "Al Wasl"
$ws.Range("I54").Value = 4
$ws.Range("J54").Value = 2.59
$ws.Range("K54").Value = "19/11/2023 21:12"
$ws.Range("L54").Value = 2.37
$ws.Range("M54").Value = "25/11/2023 16:26"
$ws.Range("N54").Value = 3.87
$ws.Range("O54").Value = "19/11/2023 21:12"
$ws.Range("P54").Value = 4.09
$ws.Range("Q54").Value = "25/11/2023 16:29"
$ws.Range("R54").Value = 2.34
$ws.Range("S54").Value = "19/11/2023 21:12"
$ws.Range("T54").Value = 2.65
$ws.Range("U54").Value = "25/11/2023 16:26"
$ws.Range("V54").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-jazira-al-wasl/bwgnEsYq/"

# Row 55 (Indice=54)
$ws.Range("B55").Value = "united-arab-emirates"
$ws.Range("C55").Value = "uae-league"
$ws.Range("D55").Value = "2023-2024"
$ws.Range("A55").Value = 54
$ws.Range("E55").Value = 45256.57291666666
$ws.Range("F55").Value = "Khorfakkan"
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = "Ajman"
$ws.Range("I55").Value = 2
$ws.Range("J55").Value = 2.39
$ws.Range("K55").Value = "20/11/2023 15:12"
$ws.Range("L55").Value = 2.18
$ws.Range("M55").Value = "26/11/2023 13:41"
$ws.Range("N55").Value = 3.68
$ws.Range("O55").Value = "20/11/2023 15:12"
$ws.Range("P55").Value = 3.6
$ws.Range("Q55").Value = "26/11/2023 13:41"
$ws.Range("R55").Value = 2.75
$ws.Range("S55").Value = "20/11/2023 15:12"
$ws.Range("T55").Value = 3.24
$ws.Range("U55").Value = "26/11/2023 13:35"
$ws.Range("V55").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/khorfakkan-ajman/dU6KINPS/"

# Row 56 (Indice=55)
$ws.Range("B56").Value = "united-arab-emirates"
$ws.Range("C56").Value = "uae-league"
$ws.Range("D56").Value = "2023-2024"
$ws.Range("A56").Value = 55
$ws.Range("E56").Value = 45256.57291666666
$ws.Range("F56").Value = "Ittihad Kalba"
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = "Al Nasr"
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 2.24
$ws.Range("K56").Value = "20/11/2023 15:12"
$ws.Range("L56").Value = 2.82
$ws.Range("M56").Value = "26/11/2023 13:35"
$ws.Range("N56").Value = 3.71
$ws.Range("O56").Value = "20/11/2023 15:12"
$ws.Range("P56").Value = 3.76
$ws.Range("Q56").Value = "26/11/2023 13:35"
$ws.Range("R56").Value = 2.81
$ws.Range("S56").Value = "20/11/2023 15:12"
$ws.Range("T56").Value = 2.36
$ws.Range("U56").Value = "26/11/2023 13:35"
$ws.Range("V56").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/ittihad-kalba-al-nasr/E5HFJ3AM/"

# Row 57 (Indice=56)
$ws.Range("B57").Value = "united-arab-emirates"
$ws.Range("C57").Value = "uae-league"
$ws.Range("D57").Value = "2023-2024"
$ws.Range("A57").Value = 56
$ws.Range("E57").Value = 45256.6875
$ws.Range("F57").Value = "Shabab Al-Ahli Dubai"
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = "Al Ain"
$ws.Range("I57").Value = 3
$ws.Range("J57").Value = 2.5
$ws.Range("K57").Value = "20/11/2023 15:12"
$ws.Range("L57").Value = 2.78
$ws.Range("M57").Value = "26/11/2023 16:21"
$ws.Range("N57").Value = 3.58
$ws.Range("O57").Value = "20/11/2023 15:12"
$ws.Range("P57").Value = 3.92
$ws.Range("Q57").Value = "26/11/2023 16:21"
$ws.Range("R57").Value = 2.67
$ws.Range("S57").Value = "20/11/2023 15:12"
$ws.Range("T57").Value = 2.33
$ws.Range("U57").Value = "26/11/2023 16:21"
$ws.Range("V57").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/shabab-al-ahli-dubai-al-ain/SCcjD1mj/"

